# Add a new "Number of cyclists" column to the Amsterdam results sheet.
# The new column is inserted right after "Area (km2)" (column D) and before
# "Beardwood approx" (old column E), which pushes every column from the old
# E onward one slot to the right (old E..M -> new F..N).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column at E - this shifts the existing E:M columns to F:N
# and keeps all their values/styles/formatting intact.
$ws.Columns.Item(5).Insert()

# New header for the inserted column.
$ws.Range("E1").Value = "Number of cyclists"

# Per-postcode cyclist counts, in row order (row 2 .. row 21).
$cyclists = @(48, 40, 1, 71, 1, 1, 104, 118, 102, 61, 1, 1, 61, 97, 80, 1, 73, 38, 41, 55)

for ($i = 0; $i -lt $cyclists.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 5).Value = $cyclists[$i]
}
